$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the three obsolete header rows (1-3) which contained the
# Sekretariat/KMK letterhead text. Everything below shifts up by 3, so the
# "Land / Sommer / Herbst / ..." table header (old row 4) becomes row 1,
# Baden-Wuerttemberg (old row 5) becomes row 2, and so on through
# Schleswig-Holstein, which moves from old row 16 to new row 13. The three
# footnote rows move from old rows 17-19 to new rows 14-16.
$ws.Rows("1:3").Delete()

# The first footnote ("*) Am 02. und 03.11.1987 ist schulfrei") used to be
# merged across A:G (old A17:G17). After the shift it lives in row 14 and
# should no longer be merged, while the remaining two footnotes (now rows
# 15 and 16) stay merged across A:G.
$ws.Range("A14:G14").UnMerge()

# Restore the original per-row heights (row height is a row-level property
# that stays keyed to its row number rather than travelling with the cut
# content), including the three now-empty trailing rows 17-19 which keep
# their old heights even though they no longer hold any text.
$ws.Rows(1).RowHeight = 27.95
$ws.Rows(2).RowHeight = 27
$ws.Rows(3).RowHeight = 18
$ws.Rows(4).RowHeight = 32.1
$ws.Rows(5).RowHeight = 21.95
$ws.Rows(6).RowHeight = 24
$ws.Rows(7).RowHeight = 21.95
$ws.Rows(8).RowHeight = 21.95
$ws.Rows(9).RowHeight = 27.95
$ws.Rows(10).RowHeight = 17.1
$ws.Rows(11).RowHeight = 21.95
$ws.Rows(12).RowHeight = 21.95
$ws.Rows(13).RowHeight = 21.95
$ws.Rows(14).RowHeight = 21.95
$ws.Rows(15).RowHeight = 24.95
$ws.Rows(16).RowHeight = 18.95
$ws.Rows(17).RowHeight = 15.95
$ws.Rows(18).RowHeight = 15
$ws.Rows(19).RowHeight = 15
